# Update the description text for the "RECURSOS RECEBIDOS PARA LIVRE UTILIZAÇÃO"
# row to the new wording "RECURSOS RECEBIDOS PARA EXECUÇÃO DIRETA DAS UNIDADES
# ORÇAMENTÁRIAS" on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "RECURSOS RECEBIDOS PARA EXECUÇÃO DIRETA DAS UNIDADES ORÇAMENTÁRIAS"

# Move/restore the active cell selection to match the saved view state.
$ws.Range("D17").Select()
